$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74..121 down to 75..122
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new weekly record
$ws.Range("A74").Value2 = 7
$ws.Range("B74").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C74").Value2 = "Ñuble"
$ws.Range("D74").Value2 = 45001
$ws.Range("E74").Value2 = 16
$ws.Range("F74").Value2 = 100112037
$ws.Range("G74").Value2 = "Cebollín"
$ws.Range("H74").Value2 = "Sin especificar"
$ws.Range("I74").Value2 = "Primera"
$ws.Range("J74").Value2 = 60
$ws.Range("K74").Value2 = 6500
$ws.Range("L74").Value2 = 7000
$ws.Range("M74").Value2 = 6750
$ws.Range("N74").Value2 = "$/paquete 36 unidades"
$ws.Range("O74").Value2 = "Provincia de Diguillín"
$ws.Range("P74").Value2 = 188
$ws.Range("Q74").Value2 = 36
$ws.Range("R74").Value2 = "Hortaliza"
